$d = $word.ActiveDocument
$d.Content.Find.Execute("10", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11", 2)
